{"js": "// Replace the date line and each \"NNN\u00f7N=\" division expression with its\n// new value. Every old value below is unique within the document, so a\n// matchCase whole-text search+replace is safe and unambiguous.\nconst replacements = [\n  [\"2024-08-23 Friday\", \"2024-08-24 Saturday\"],\n  [\"600\u00f72=\", \"796\u00f77=\"],\n  [\"380\u00f78=\", \"821\u00f76=\"],\n  [\"463\u00f73=\", \"653\u00f75=\"],\n  [\"446\u00f78=\", \"284\u00f79=\"],\n  [\"295\u00f77=\", \"543\u00f72=\"],\n  [\"824\u00f74=\", \"311\u00f72=\"],\n  [\"406\u00f73=\", \"348\u00f78=\"],\n  [\"451\u00f72=\", \"709\u00f72=\"],\n  [\"293\u00f75=\", \"378\u00f73=\"],\n  [\"198\u00f76=\", \"682\u00f77=\"],\n  [\"431\u00f72=\", \"466\u00f78=\"],\n  [\"919\u00f78=\", \"184\u00f76=\"],\n  [\"237\u00f75=\", \"667\u00f74=\"],\n  [\"558\u00f73=\", \"291\u00f79=\"],\n  [\"851\u00f79=\", \"968\u00f74=\"],\n  [\"376\u00f76=\", \"189\u00f75=\"],\n  [\"968\u00f76=\", \"165\u00f72=\"],\n  [\"164\u00f74=\", \"989\u00f76=\"],\n  [\"849\u00f73=\", \"772\u00f72=\"],\n  [\"900\u00f78=\", \"914\u00f76=\"],\n  [\"509\u00f77=\", \"608\u00f75=\"],\n  [\"526\u00f74=\", \"827\u00f76=\"],\n  [\"727\u00f78=\", \"688\u00f78=\"],\n  [\"981\u00f76=\", \"293\u00f77=\"],\n  [\"418\u00f77=\", \"913\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00f7N=\" division expression with its\n# new value. Every old value below is unique within the document, so a\n# MatchCase Find/Replace is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-08-23 Friday\", \"2024-08-24 Saturday\"),\n  @(\"600\u00f72=\", \"796\u00f77=\"),\n  @(\"380\u00f78=\", \"821\u00f76=\"),\n  @(\"463\u00f73=\", \"653\u00f75=\"),\n  @(\"446\u00f78=\", \"284\u00f79=\"),\n  @(\"295\u00f77=\", \"543\u00f72=\"),\n  @(\"824\u00f74=\", \"311\u00f72=\"),\n  @(\"406\u00f73=\", \"348\u00f78=\"),\n  @(\"451\u00f72=\", \"709\u00f72=\"),\n  @(\"293\u00f75=\", \"378\u00f73=\"),\n  @(\"198\u00f76=\", \"682\u00f77=\"),\n  @(\"431\u00f72=\", \"466\u00f78=\"),\n  @(\"919\u00f78=\", \"184\u00f76=\"),\n  @(\"237\u00f75=\", \"667\u00f74=\"),\n  @(\"558\u00f73=\", \"291\u00f79=\"),\n  @(\"851\u00f79=\", \"968\u00f74=\"),\n  @(\"376\u00f76=\", \"189\u00f75=\"),\n  @(\"968\u00f76=\", \"165\u00f72=\"),\n  @(\"164\u00f74=\", \"989\u00f76=\"),\n  @(\"849\u00f73=\", \"772\u00f72=\"),\n  @(\"900\u00f78=\", \"914\u00f76=\"),\n  @(\"509\u00f77=\", \"608\u00f75=\"),\n  @(\"526\u00f74=\", \"827\u00f76=\"),\n  @(\"727\u00f78=\", \"688\u00f78=\"),\n  @(\"981\u00f76=\", \"293\u00f77=\"),\n  @(\"418\u00f77=\", \"913\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
